# edit.ps1 - apply the changes described by the diff:
#  1) Slide 6's table switches to a different (built-in) table style GUID.
#  2) The presentation's active theme (the one actually driving the
#     slide master / slides, i.e. ppt/theme/theme2.xml) swaps its
#     12-colour palette from the "Integral" scheme to the classic
#     "Office" scheme, so visually the deck goes from the green/teal
#     Integral look back to the default Office blue/orange look.

$p = $ppt.ActivePresentation

# --- 1) Update the table style on slide 6 -----------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{29696537-3172-4856-B5AA-1D53DF2A7635}")
    }
}

# --- 2) Swap the theme colour palette back to "Office" ----------------
function HexToRgbInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Order matches ThemeColorScheme.Colors(1..12):
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $colorScheme.Colors($i).RGB = HexToRgbInt($officeColors[$i - 1])
}
